# 15/02/2015 : Reconnexion avec la base de donnees
# - Rewrite the "Donnees" reference-data rows (DemandeStatus / VetementType /
#   CompteurType / NotifType insert statements) with the corrected SQL.
# - Rename "Sheet3" to "Identifiants de test" and give it test credentials.
# - Add a new "Script divers" sheet with a DBCC CHECKIDENT repair script for
#   the Personne table identity-seed bug mentioned in the commit message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Donnees" sheet: replace the TypeParam reference rows (35-61) so the
#    DemandeStatus block uses dbo.TypeParam + longer french labels, and the
#    VetementType / CompteurType / NotifType blocks are reordered/extended.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Données")

$ws2.Range("B35").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AttenteDemandeur','Proposition en attente réponse demandeur')"
$ws2.Range("B36").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AttenteConseiller','Demande en attente réponse conseiller')"
$ws2.Range("B37").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','Accepte','Demande ou proposition acceptée')"
$ws2.Range("B38").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','RefusDemandeur','Proposition refusée par demandeur')"
$ws2.Range("B39").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','RefusConseiller','Demande refusée par conseiller')"
$ws2.Range("B40").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AnnulAdmin','Demande ou proposition annulée par administrateur')"
$ws2.Range("B41").Value = "insert into dbo.TypeParam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','Termine','Demande ou proposition terminée')"

$ws2.Range("B43").Value = "Type de vêtement"
$ws2.Range("E43").Value = "VetementType"
$ws2.Range("B44").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Tete','Tête')"
$ws2.Range("B45").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Buste','Buste')"
$ws2.Range("B46").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Jambe','Jambe')"
$ws2.Range("B47").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Pied','Pied')"
$ws2.Range("B48").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Accessoire','Accessoire')"
$ws2.Range("B49").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (5,'VetementType','Main','Main')"

# Old B50/E50 ("Type compteur" header) moves down one row to B51/E51 - clear
# the old cells (including their bold/underline header style) first.
$ws2.Range("B50").Clear()
$ws2.Range("E50").Clear()

$ws2.Range("B51").Value = "Type compteur"
$ws2.Range("E51").Value = "CompteurType"
$ws2.Range("B51").Font.Bold = $true
$ws2.Range("B51").Font.Underline = 2
$ws2.Range("B52").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (6,'CompteurType','Abonne','Abonné')"
$ws2.Range("B53").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (6,'CompteurType','Conseiller','Conseiller')"

# Old B54/E54 ("Type de notification" header) moves down to B55/E55 - clear
# the old cells first.
$ws2.Range("B54").Clear()
$ws2.Range("E54").Clear()

$ws2.Range("B55").Value = "Type de notification"
$ws2.Range("E55").Value = "NotifType"
$ws2.Range("B55").Font.Bold = $true
$ws2.Range("B55").Font.Underline = 2
$ws2.Range("B56").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','DemandeCreation','Demander une aide')"
$ws2.Range("B57").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','DemandeAccept','Demande acceptée')"
$ws2.Range("B58").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','DemandeReject','Demande rejetée')"
$ws2.Range("B59").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','PropositionAccept','Proposition acceptée')"
$ws2.Range("B60").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','PropositionReject','Proposition rejetée')"
$ws2.Range("B61").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (7,'NotifType','PropositionCreation','Proposer une aide')"

# ---------------------------------------------------------------------------
# 2) Rename "Sheet3" -> "Identifiants de test" and fill in test credentials.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "Identifiants de test"
$ws3.Range("A1").Value = "userA"
$ws3.Range("B1").Value = "mdpmdpA"

# ---------------------------------------------------------------------------
# 3) Add the new "Script divers" sheet (after "Identifiants de test") with a
#    DBCC CHECKIDENT repair script for the Personne identity-seed bug.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Script divers"

$ws4.Range("A2").Value = "DBCC CHECKIDENT ('dbo.UserProfile');"
$ws4.Range("A3").Value = "GO"
$ws4.Range("A4").Value = "DBCC CHECKIDENT ('dbo.webpages_Membership');"
$ws4.Range("A5").Value = "GO"
$ws4.Range("A6").Value = "USE ModeConseil"
$ws4.Range("A7").Value = "GO"
$ws4.Range("A8").Value = "DBCC CHECKIDENT ('Personne');"
$ws4.Range("A9").Value = "GO"

# ---------------------------------------------------------------------------
# 4) Selections per sheet + which tab ends up active (matches the saved view
#    state captured in the workbook).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B56:B61").Select()

$ws3.Activate()
$ws3.Range("B2").Select()

$ws4.Activate()
$ws4.Range("A2:A9").Select()
